$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with rows 194-196 holding the most recent week's
# Esparragos quotes (date 44477). A new week's data (date 45239) is being
# added "on top" of rows 194-196, and the previous week's rows get pushed
# down to become new rows 197-199 (unchanged copies of the old 194-196).

# 1) Duplicate the existing rows 194-196 down to 197-199 (exact copies,
#    preserving old date/values) before overwriting 194-196 with the new week.
$ws.Rows.Item(194).Copy()
$ws.Rows.Item(197).Insert()

$ws.Rows.Item(195).Copy()
$ws.Rows.Item(198).Insert()

$ws.Rows.Item(196).Copy()
$ws.Rows.Item(199).Insert()

# 2) Overwrite rows 194-196 (the "Banquete"/"Primera"/"Segunda" quality rows)
#    with the new week's date and updated figures.

# Row 194 - Banquete: date + price update (volume unchanged)
$ws.Range("D194").Value = 45239
$ws.Range("K194").Value = 1700
$ws.Range("L194").Value = 1700
$ws.Range("M194").Value = 1700
$ws.Range("P194").Value = 1700

# Row 195 - Primera: date + volume + price update
$ws.Range("D195").Value = 45239
$ws.Range("J195").Value = 520
$ws.Range("K195").Value = 1500
$ws.Range("L195").Value = 1500
$ws.Range("M195").Value = 1500
$ws.Range("P195").Value = 1500

# Row 196 - Segunda: date + volume + price update
$ws.Range("D196").Value = 45239
$ws.Range("J196").Value = 340
$ws.Range("K196").Value = 1300
$ws.Range("L196").Value = 1300
$ws.Range("M196").Value = 1300
$ws.Range("P196").Value = 1300
